$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Jacques Friesen"
$ws.Range("C2").Value = "lavern.ryan@example.com"
$ws.Range("D2").Value = "2025-06-22T19:13:32.000000Z"
$ws.Range("E2").Value = "2025-06-22T19:13:32.000000Z"
$ws.Range("F2").Value = "2025-06-22T19:13:32.000000Z"

# Row 3
$ws.Range("B3").Value = "Marilie Wiegand"
$ws.Range("C3").Value = "feest.quinten@example.net"
$ws.Range("D3").Value = "2025-06-22T19:13:32.000000Z"
$ws.Range("E3").Value = "2025-06-22T19:13:32.000000Z"
$ws.Range("F3").Value = "2025-06-22T19:13:32.000000Z"

# Row 4
$ws.Range("B4").Value = "Dr. Nia Kutch"
$ws.Range("C4").Value = "amanda03@example.net"
$ws.Range("D4").Value = "2025-06-22T19:13:32.000000Z"
$ws.Range("E4").Value = "2025-06-22T19:13:32.000000Z"
$ws.Range("F4").Value = "2025-06-22T19:13:32.000000Z"
